$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 6

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 6

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 11
$ws.Range("D4").Value = 13

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = 19

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = 25

# Row 7 (new row)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 31
